$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-11, columns B-G
$data = @(
    @{Row=2;  B=0.1155809658704365;  C=0.390457646890028;   D=0.3755683084142444; E=0.612836281901002;  F=0.6153644530269748; G=23}
    @{Row=3;  B=0.6511545689434494;  C=0.8016979307476643;  D=4.674211074126436;  E=2.161992385307228;  F=2.107937406714976;  G=23}
    @{Row=4;  B=0.2201423283843545;  C=1.112030744985122;   D=7.605937587058003;  E=2.757886434764492;  F=2.810871133758786;  G=23}
    @{Row=5;  B=0.2501754929190889;  C=1.125383142428414;   D=7.653888439551237;  E=2.766566182029853;  F=2.817154554387687;  G=23}
    @{Row=6;  B=0.3004374354251687;  C=1.183046795173429;   D=7.775752270546474;  E=2.788503589839266;  F=2.834577579661701;  G=23}
    @{Row=7;  B=0.2656501924638524;  C=1.272698581772992;   D=7.926928235987235;  E=2.815480107546;     F=2.86591440951239;   G=23}
    @{Row=8;  B=0.1735487773573387; C=1.27905897888513;    D=7.992029374264993;  E=2.827017752732549; F=2.885102286758365; G=23}
    @{Row=9;  B=0.2219982892021358;  C=1.341063655789804;   D=8.0376137195044;    E=2.835068556402896;  F=2.889885189754726;  G=23}
    @{Row=10; B=0.1895698287473124;  C=1.293626197483444;   D=8.015968031220561;  E=2.831248493371884;  F=2.88838363227864;   G=23}
    @{Row=11; B=0.08322311406686993; C=1.20051888709089;    D=7.697829717014696;  E=2.77449629969382;   F=2.835575807300819;  G=23}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("B$r").Value = $rowData.B
    $ws.Range("C$r").Value = $rowData.C
    $ws.Range("D$r").Value = $rowData.D
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("F$r").Value = $rowData.F
    $ws.Range("G$r").Value = $rowData.G
}
